# Update the data row (row 2) of the extract-results sheet with a fresh
# sample record. Secondary-zip-code handling now continues on error, and
# the record below reflects the corrected/current extraction for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold values which look like numbers/dates (CollectionDate,
# RequisitionNumber, PMS, EMR, DateOfBirth, ZipCode, HomePhone) get a
# leading apostrophe so Excel keeps them as literal text, the same way
# the source extract stores them (as shared strings), instead of
# auto-converting to a numeric/date cell value.
$ws.Range("B2").Value = "'2024-02-14"
$ws.Range("C2").Value = "'271436"
$ws.Range("D2").Value = "'13441336"
$ws.Range("E2").Value = "'921624"
$ws.Range("G2").Value = "BERNSTEIN"
$ws.Range("H2").Value = "MELVIN"
$ws.Range("J2").Value = "'1947-06-08"
$ws.Range("K2").Value = "Male"
$ws.Range("L2").Value = "1113 HAMPSTEAD LN"
$ws.Range("N2").Value = "ALLEN"
$ws.Range("O2").Value = "'75013"
$ws.Range("P2").Value = "'8327228681"
$ws.Range("R2").Value = "Christopher Stroud, MD,"
$ws.Range("S2").Value = "(MK) Dallas Assoc Derm"
$ws.Range("W2").Value = "MEDICARE PART B"
